$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1832
$ws.Range("J17").Value = 1853.3334
$ws.Range("L17").Value = 5560.0002
$ws.Range("N17").Value = -5896.0002
$ws.Range("H70").Value = 1077.1613
$ws.Range("I70").Value = 849.61536
$ws.Range("J70").Value = 1241.5
$ws.Range("K70").Value = 2548.84608
$ws.Range("L70").Value = 3724.5
$ws.Range("M70").Value = -2278.84608
$ws.Range("N70").Value = -4264.5
$ws.Range("H73").Value = 1077.1613
$ws.Range("I73").Value = 849.61536
$ws.Range("J73").Value = 1241.5
$ws.Range("K73").Value = 2548.84608
$ws.Range("L73").Value = 3724.5
$ws.Range("M73").Value = -1612.84608
$ws.Range("N73").Value = -5596.5
$ws.Range("H112").Value = 1341.6666
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 1490
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 4470
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -6686
$ws.Range("H138").Value = 4001818.8
$ws.Range("I138").Value = 1334.9395
$ws.Range("J138").Value = 11767464
$ws.Range("K138").Value = 4004.8185
$ws.Range("L138").Value = 35302392
$ws.Range("M138").Value = 1135.1815
$ws.Range("N138").Value = -35312672

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7577103
$ws.Range("I45").Value = 11364654
$ws.Range("K45").Value = 11364654
$ws.Range("M45").Value = -11364277
$ws.Range("H61").Value = 2265.5881
$ws.Range("I61").Value = 1893.9286
$ws.Range("K61").Value = 1893.9286
$ws.Range("M61").Value = -1681.9286
$ws.Range("H102").Value = 3461.5
$ws.Range("I102").Value = 2998
$ws.Range("K102").Value = 2998
$ws.Range("M102").Value = -1376
$ws.Range("H136").Value = 2265.5881
$ws.Range("I136").Value = 1893.9286
$ws.Range("K136").Value = 5681.7858
$ws.Range("M136").Value = -3131.7858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2043.5172
$ws.Range("I86").Value = 1602.4445
$ws.Range("J86").Value = 2765.2727
$ws.Range("K86").Value = 1602.4445
$ws.Range("L86").Value = 2765.2727
$ws.Range("M86").Value = -479.4445000000001
$ws.Range("N86").Value = -5011.2727
$ws.Range("H89").Value = 2043.5172
$ws.Range("I89").Value = 1602.4445
$ws.Range("J89").Value = 2765.2727
$ws.Range("K89").Value = 8012.2225
$ws.Range("L89").Value = 13826.3635
$ws.Range("M89").Value = -2396.2225
$ws.Range("N89").Value = -25058.3635
$ws.Range("H94").Value = 1127
$ws.Range("I94").Value = 893.6
$ws.Range("J94").Value = 1627.1428
$ws.Range("K94").Value = 893.6
$ws.Range("L94").Value = 1627.1428
$ws.Range("M94").Value = -442.6
$ws.Range("N94").Value = -2529.1428
$ws.Range("H105").Value = 2193
$ws.Range("I105").Value = 2190.0925
$ws.Range("J105").Value = 2212.625
$ws.Range("K105").Value = 2190.0925
$ws.Range("L105").Value = 2212.625
$ws.Range("M105").Value = -443.0925000000002
$ws.Range("N105").Value = -5706.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2190.4138
$ws.Range("I31").Value = 1175.6154
$ws.Range("J31").Value = 3014.9375
$ws.Range("K31").Value = 1175.6154
$ws.Range("L31").Value = 3014.9375
$ws.Range("M31").Value = -880.6153999999999
$ws.Range("N31").Value = -3604.9375
$ws.Range("H34").Value = 2190.4138
$ws.Range("I34").Value = 1175.6154
$ws.Range("J34").Value = 3014.9375
$ws.Range("K34").Value = 1175.6154
$ws.Range("L34").Value = 3014.9375
$ws.Range("M34").Value = -973.6153999999999
$ws.Range("N34").Value = -3418.9375
$ws.Range("H99").Value = 1958.0286
$ws.Range("I99").Value = 1761.3214
$ws.Range("K99").Value = 1761.3214
$ws.Range("M99").Value = -263.3214
$ws.Range("H126").Value = 1958.0286
$ws.Range("I126").Value = 1761.3214
$ws.Range("K126").Value = 5283.9642
$ws.Range("M126").Value = -2813.9642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1489.2354
$ws.Range("I121").Value = 750
$ws.Range("J121").Value = 1587.8
$ws.Range("K121").Value = 2250
$ws.Range("L121").Value = 4763.4
$ws.Range("M121").Value = -940
$ws.Range("N121").Value = -7383.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 144.8125
$ws.Range("I2").Value = 36.2
$ws.Range("J2").Value = 194.18182
$ws.Range("K2").Value = 36.2
$ws.Range("L2").Value = 194.18182
$ws.Range("M2").Value = 76.8
$ws.Range("N2").Value = -420.18182
$ws.Range("H70").Value = 5947.2085
$ws.Range("I70").Value = 5408.8
$ws.Range("J70").Value = 14023.333
$ws.Range("K70").Value = 5408.8
$ws.Range("L70").Value = 14023.333
$ws.Range("M70").Value = -5138.8
$ws.Range("N70").Value = -14563.333
$ws.Range("H73").Value = 5947.2085
$ws.Range("I73").Value = 5408.8
$ws.Range("J73").Value = 14023.333
$ws.Range("K73").Value = 5408.8
$ws.Range("L73").Value = 14023.333
$ws.Range("M73").Value = -4472.8
$ws.Range("N73").Value = -15895.333
$ws.Range("H122").Value = 1566.2778
$ws.Range("I122").Value = 1588.5555
$ws.Range("K122").Value = 4765.666499999999
$ws.Range("M122").Value = -2315.666499999999
$ws.Range("H126").Value = 47631144
$ws.Range("I126").Value = 100006660
$ws.Range("J126").Value = 17045.455
$ws.Range("K126").Value = 300019980
$ws.Range("L126").Value = 51136.36500000001
$ws.Range("M126").Value = -300017510
$ws.Range("N126").Value = -56076.36500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1055.2941
$ws.Range("I82").Value = 985.6667
$ws.Range("J82").Value = 1070.2142
$ws.Range("K82").Value = 985.6667
$ws.Range("L82").Value = 1070.2142
$ws.Range("M82").Value = -624.6667
$ws.Range("N82").Value = -1792.2142
$ws.Range("H85").Value = 1055.2941
$ws.Range("I85").Value = 985.6667
$ws.Range("J85").Value = 1070.2142
$ws.Range("K85").Value = 985.6667
$ws.Range("L85").Value = 1070.2142
$ws.Range("M85").Value = 262.3333
$ws.Range("N85").Value = -3566.2142
$ws.Range("H100").Value = 1327.2
$ws.Range("I100").Value = 1043.2858
$ws.Range("J100").Value = 1575.625
$ws.Range("K100").Value = 1043.2858
$ws.Range("L100").Value = 1575.625
$ws.Range("M100").Value = -502.2858000000001
$ws.Range("N100").Value = -2657.625
$ws.Range("H132").Value = 3256.9387
$ws.Range("I132").Value = 3134.2307
$ws.Range("J132").Value = 3735.5
$ws.Range("K132").Value = 9402.6921
$ws.Range("L132").Value = 11206.5
$ws.Range("M132").Value = -6872.6921
$ws.Range("N132").Value = -16266.5
$ws.Range("H136").Value = 2049.4
$ws.Range("I136").Value = 1783.4286
$ws.Range("K136").Value = 5350.2858
$ws.Range("M136").Value = -2800.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 43736.32
$ws.Range("I122").Value = 770.1739
$ws.Range("K122").Value = 2310.5217
$ws.Range("M122").Value = 139.4782999999998
